$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + week-of date range) ---
$ws.Range("A8").Value = "Volume 32   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/10/2025  Through  2/16/2025"

# --- Cells whose style/type flips between numeric and text placeholder ---
# Donor cells (style never changes across the diff): C17 (numeric #,##0 style),
# D14 (text "0" placeholder style), E14 (text "***.*" placeholder style)
$ws.Range("C17").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("C17").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 2
$ws.Range("D14").Copy($ws.Range("C18"))
$ws.Range("D14").Copy($ws.Range("C20"))
$ws.Range("D14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("D14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("C17").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1

# --- Plain value updates (style unchanged) ---
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 200
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = -50
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 26
$ws.Range("K16").Value = -57.692307692307
$ws.Range("L16").Value = -35.294117647058
$ws.Range("M16").Value = -67.647058823529
$ws.Range("N16").Value = -86.25
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -10
$ws.Range("I17").Value = 17
$ws.Range("J17").Value = 22
$ws.Range("K17").Value = -22.727272727272
$ws.Range("L17").Value = -22.727272727272
$ws.Range("M17").Value = 30.769230769230
$ws.Range("N17").Value = -56.410256410256
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 1
$ws.Range("H18").Value = -85.714285714285
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = -54.545454545454
$ws.Range("L18").Value = -68.75
$ws.Range("M18").Value = -81.481481481481
$ws.Range("N18").Value = -96.183206106870
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 28
$ws.Range("H19").Value = -20
$ws.Range("I19").Value = 63
$ws.Range("J19").Value = 57
$ws.Range("K19").Value = 10.526315789473
$ws.Range("L19").Value = -8.695652173913
$ws.Range("M19").Value = 3.278688524590
$ws.Range("N19").Value = -19.230769230769
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("M20").Value = -34.782608695652
$ws.Range("N20").Value = -96.287128712871
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -26.315789473684
$ws.Range("F21").Value = 53
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -28.378378378378
$ws.Range("I21").Value = 114
$ws.Range("J21").Value = 130
$ws.Range("K21").Value = -12.307692307692
$ws.Range("L21").Value = -14.925373134328
$ws.Range("M21").Value = -28.301886792452
$ws.Range("N21").Value = -84.594594594594
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -29.729729729729
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 139
$ws.Range("H24").Value = -19.424460431654
$ws.Range("I24").Value = 183
$ws.Range("J24").Value = 188
$ws.Range("K24").Value = -2.659574468085
$ws.Range("L24").Value = 40.769230769230
$ws.Range("M24").Value = 57.758620689655
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -38.461538461538
$ws.Range("F25").Value = 70
$ws.Range("G25").Value = 99
$ws.Range("H25").Value = -29.292929292929
$ws.Range("I25").Value = 119
$ws.Range("J25").Value = 135
$ws.Range("K25").Value = -11.851851851851
$ws.Range("L25").Value = 33.707865168539
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -16.666666666666
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -46.428571428571
$ws.Range("I26").Value = 29
$ws.Range("J26").Value = 42
$ws.Range("K26").Value = -30.952380952381
$ws.Range("L26").Value = -23.684210526315
$ws.Range("M26").Value = -23.684210526315
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = -25
$ws.Range("L29").Value = -75
$ws.Range("N29").Value = -85.714285714285
$ws.Range("L30").Value = -66.666666666666
$ws.Range("N30").Value = -85.714285714285
